$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round row 5 (B5:AH5) values to 2 decimal places using "round half away from zero" (Excel ROUND semantics)
for ($col = 2; $col -le 34; $col++) {
    $cell = $ws.Cells.Item(5, $col)
    $val = $cell.Value()
    if ($val -ge 0) {
        $rounded = [Math]::Floor(($val * 100) + 0.5) / 100
    } else {
        $rounded = [Math]::Ceiling(($val * 100) - 0.5) / 100
    }
    $cell.Value = $rounded
}

# Delete entire row 6
$ws.Rows.Item(6).Delete()
